# Generate Report for Handback
#
# Row 7 on the "zh-cn" and "de-de" sheets corresponds to the handback of
# 3ec2fd86-21b5-4a8e-887e-83b30033a907.md. A new handback was processed,
# but its committed version was stale, so:
#   - I7 (Latest Target File) now links to the handback .md file itself
#   - J7 (Latest Handback File) is populated with the generated .xlf name
#   - K7 (Latest Handback DateTime) gets the handback-processed timestamp
#   - P7 (Error Detail) explains the stale-version problem

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$mdName = "3ec2fd86-21b5-4a8e-887e-83b30033a907.md"
$mdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb2a67bdfb16170033929aad4cb26ce4fcb05d85/e2e/3ec2fd86-21b5-4a8e-887e-83b30033a907.md"
$errMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3d857b56135beff173024331b8a326a1a66bbab/e2e/3ec2fd86-21b5-4a8e-887e-83b30033a907.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb2a67bdfb16170033929aad4cb26ce4fcb05d85/e2e/3ec2fd86-21b5-4a8e-887e-83b30033a907.md."

# Same look as the workbook's other handoff/handback hyperlinks (custom
# "HyperLink" cell style: single underline, cornflower-blue font color).
$linkUnderline = 2         # xlUnderlineStyleSingle
$linkColor     = 15570276  # RGB(100,149,237) == #6495ED, BGR-packed long

# ---------------- zh-cn sheet, row 7 ----------------
$zhcn.Range("I7").Value = $mdName
$zhcn.Hyperlinks.Add($zhcn.Range("I7"), $mdUrl, "", "", $mdName) | Out-Null
# Hyperlinks.Add stamps its own (theme-based) link formatting, so apply the
# workbook's custom link color/underline afterwards.
$zhcn.Range("I7").Font.Underline = $linkUnderline
$zhcn.Range("I7").Font.Color = $linkColor

$zhcn.Range("J7").Value = "3ec2fd86-21b5-4a8e-887e-83b30033a907.5f0530ed4f9bba00c0051961585552ed3bb74f28.zh-cn.xlf"
$zhcn.Range("K7").Value = "2016-08-30 17:04:54"
$zhcn.Range("P7").Value = $errMsg

# ---------------- de-de sheet, row 7 ----------------
$dede.Range("I7").Value = $mdName
$dede.Hyperlinks.Add($dede.Range("I7"), $mdUrl, "", "", $mdName) | Out-Null
$dede.Range("I7").Font.Underline = $linkUnderline
$dede.Range("I7").Font.Color = $linkColor

$dede.Range("J7").Value = "3ec2fd86-21b5-4a8e-887e-83b30033a907.5f0530ed4f9bba00c0051961585552ed3bb74f28.de-de.xlf"
$dede.Range("K7").Value = "2016-08-30 17:05:08"
$dede.Range("P7").Value = $errMsg
